# Update the "cryptos" listing on Sheet1 with refreshed price/volume data
# (and a refreshed set of bottom-ranked coins), as produced by the
# scheduled GitHub Actions refresh job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds text-formatted numbers (e.g. "584.73",
# "2.569.23" using '.' as both thousands- and decimal-separator look-alikes).
# Force the cell format to Text ("@") before assigning so Excel does not
# silently reinterpret the string as a number (which would also normalise
# away things like trailing zeros or turn "0.0000147" into scientific
# notation).

# --- Price (column D) / Volume(1h) (column E) refreshes ---------------
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.246.47"
$ws.Range("E2").Value = "  +0.79%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.569.59"
$ws.Range("E3").Value = "  +1.15%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.64"
$ws.Range("E5").Value = "  +3.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.13"
$ws.Range("E6").Value = "  +1.63%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.601"
$ws.Range("E8").Value = "  +3.60%  "
$ws.Range("E9").Value = "  +4.00%  "
$ws.Range("E10").Value = "  +0.81%  "
$ws.Range("E11").Value = "  +0.47%  "
$ws.Range("E12").Value = "  +1.52%  "
$ws.Range("E13").Value = "  +1.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.027.63"
$ws.Range("E14").Value = "  +1.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.188.84"
$ws.Range("E15").Value = "  +0.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000147"
$ws.Range("E16").Value = "  +4.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.553.41"
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("E18").Value = "  -0.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "343.04"
$ws.Range("E19").Value = "  +2.85%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.41"
$ws.Range("E20").Value = "  +3.44%  "
$ws.Range("E21").Value = "  +2.18%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.70"
$ws.Range("E23").Value = "  +3.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.694.14"
$ws.Range("E24").Value = "  +1.58%  "
$ws.Range("E25").Value = "  +3.29%  "
$ws.Range("E26").Value = "  +1.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.20"
$ws.Range("E27").Value = "  +13.22%  "
$ws.Range("E28").Value = "  +2.19%  "
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("E31").Value = "  +7.92%  "
$ws.Range("E32").Value = "  +2.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "461.01"
$ws.Range("E34").Value = "  +4.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "176.63"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.406"
$ws.Range("E36").Value = "  +2.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.24"
$ws.Range("E37").Value = "  +1.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.53"
$ws.Range("E38").Value = "  +4.49%  "
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "151.41"
$ws.Range("E42").Value = "  -0.84%  "
$ws.Range("E43").Value = "  +2.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.08"
$ws.Range("E44").Value = "  +2.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0553"
$ws.Range("E45").Value = "  +7.07%  "
$ws.Range("E46").Value = "  +1.98%  "
$ws.Range("E47").Value = "  +2.78%  "
$ws.Range("E48").Value = "  +2.27%  "

# --- Rows 49-51: EnergySwap dropped out of the ranking; the remaining
#     coins shift up one slot and TheGraph enters at the bottom -------
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.76"
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.38"
$ws.Range("E50").Value = "  -0.12%  "
$ws.Range("B51").Value = "TheGraph"
$ws.Range("C51").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.164"
$ws.Range("E51").Value = "  +4.47%  "
